# Auto-applied scheduled market-data refresh for the Leve profit workbook.
# For each affected sheet/row, columns H:N (currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ), LeveProfit(NQ/HQ)) are overwritten with refreshed values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 72.77778000000001
$ws.Range("I11").Value = 72.77778000000001
$ws.Range("K11").Value = 72.77778000000001
$ws.Range("M11").Value = 67.22221999999999
$ws.Range("H19").Value = 1135.5834
$ws.Range("I19").Value = 749.3333
$ws.Range("J19").Value = 1264.3334
$ws.Range("K19").Value = 749.3333
$ws.Range("L19").Value = 1264.3334
$ws.Range("M19").Value = -574.3333
$ws.Range("N19").Value = -1614.3334
$ws.Range("H28").Value = 6170.4287
$ws.Range("I28").Value = 1037.2
$ws.Range("K28").Value = 1037.2
$ws.Range("M28").Value = -552.2
$ws.Range("H40").Value = 9374.651
$ws.Range("I40").Value = 7688.1333
$ws.Range("J40").Value = 10278.143
$ws.Range("K40").Value = 7688.1333
$ws.Range("L40").Value = 10278.143
$ws.Range("M40").Value = -7513.1333
$ws.Range("N40").Value = -10628.143
$ws.Range("H138").Value = 3973.2354
$ws.Range("I138").Value = 2354.5557
$ws.Range("K138").Value = 7063.6671
$ws.Range("M138").Value = -1923.6671
$ws.Range("H141").Value = 4893.3887
$ws.Range("I141").Value = 3009.111
$ws.Range("K141").Value = 9027.332999999999
$ws.Range("M141").Value = -3847.332999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4538.294
$ws.Range("I63").Value = 1378.7142
$ws.Range("J63").Value = 6750
$ws.Range("K63").Value = 1378.7142
$ws.Range("L63").Value = 6750
$ws.Range("M63").Value = -692.7141999999999
$ws.Range("N63").Value = -8122
$ws.Range("H66").Value = 4538.294
$ws.Range("I66").Value = 1378.7142
$ws.Range("J66").Value = 6750
$ws.Range("K66").Value = 6893.571
$ws.Range("L66").Value = 33750
$ws.Range("M66").Value = -3461.571
$ws.Range("N66").Value = -40614
$ws.Range("H74").Value = 10419904
$ws.Range("I74").Value = 10755417
$ws.Range("J74").Value = 19014
$ws.Range("K74").Value = 10755417
$ws.Range("L74").Value = 19014
$ws.Range("M74").Value = -10754543
$ws.Range("N74").Value = -20762
$ws.Range("H77").Value = 10419904
$ws.Range("I77").Value = 10755417
$ws.Range("J77").Value = 19014
$ws.Range("K77").Value = 53777085
$ws.Range("L77").Value = 95070
$ws.Range("M77").Value = -53772717
$ws.Range("N77").Value = -103806
$ws.Range("H122").Value = 3252.5173
$ws.Range("I122").Value = 2122.6316
$ws.Range("K122").Value = 6367.8948
$ws.Range("M122").Value = -3917.8948

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1800.75
$ws.Range("I5").Value = 1052
$ws.Range("K5").Value = 1052
$ws.Range("M5").Value = -939
$ws.Range("H99").Value = 3310.2
$ws.Range("I99").Value = 3233.5557
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 3233.5557
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -1735.5557
$ws.Range("N99").Value = -6996
$ws.Range("H105").Value = 11428.821
$ws.Range("I105").Value = 11128.542
$ws.Range("K105").Value = 11128.542
$ws.Range("M105").Value = -9381.541999999999
$ws.Range("H131").Value = 55000
$ws.Range("J131").Value = 55000
$ws.Range("L131").Value = 55000
$ws.Range("N131").Value = -65080
$ws.Range("H134").Value = 2198.2449
$ws.Range("I134").Value = 1702.1957
$ws.Range("J134").Value = 9804.333000000001
$ws.Range("K134").Value = 5106.5871
$ws.Range("L134").Value = 29412.999
$ws.Range("M134").Value = -2571.5871
$ws.Range("N134").Value = -34482.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29070.637
$ws.Range("J31").Value = 51119.125
$ws.Range("L31").Value = 51119.125
$ws.Range("N31").Value = -51709.125
$ws.Range("H34").Value = 29070.637
$ws.Range("J34").Value = 51119.125
$ws.Range("L34").Value = 51119.125
$ws.Range("N34").Value = -51523.125
$ws.Range("H62").Value = 13183.846
$ws.Range("J62").Value = 14001
$ws.Range("L62").Value = 14001
$ws.Range("N62").Value = -15249
$ws.Range("H65").Value = 13183.846
$ws.Range("J65").Value = 14001
$ws.Range("L65").Value = 70005
$ws.Range("N65").Value = -76245
$ws.Range("H122").Value = 4007.4424
$ws.Range("I122").Value = 2377.1082
$ws.Range("J122").Value = 8028.933
$ws.Range("K122").Value = 7131.3246
$ws.Range("L122").Value = 24086.799
$ws.Range("M122").Value = -4681.3246
$ws.Range("N122").Value = -28986.799

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3341456.8
$ws.Range("I4").Value = 4538856.5
$ws.Range("J4").Value = 1246007.5
$ws.Range("K4").Value = 13616569.5
$ws.Range("L4").Value = 3738022.5
$ws.Range("M4").Value = -13616457.5
$ws.Range("N4").Value = -3738246.5
$ws.Range("H9").Value = 161517
$ws.Range("J9").Value = 3820.4
$ws.Range("L9").Value = 11461.2
$ws.Range("N9").Value = -11909.2
$ws.Range("H64").Value = 111116410
$ws.Range("I64").Value = 142861020
$ws.Range("J64").Value = 10257
$ws.Range("K64").Value = 428583060
$ws.Range("L64").Value = 30771
$ws.Range("M64").Value = -428582790
$ws.Range("N64").Value = -31311
$ws.Range("H67").Value = 111116410
$ws.Range("I67").Value = 142861020
$ws.Range("J67").Value = 10257
$ws.Range("K67").Value = 428583060
$ws.Range("L67").Value = 30771
$ws.Range("M67").Value = -428582124
$ws.Range("N67").Value = -32643
$ws.Range("H68").Value = 5190.077
$ws.Range("J68").Value = 5754.2173
$ws.Range("L68").Value = 17262.6519
$ws.Range("N68").Value = -18884.6519
$ws.Range("H71").Value = 5190.077
$ws.Range("J71").Value = 5754.2173
$ws.Range("L71").Value = 51787.95570000001
$ws.Range("N71").Value = -59899.95570000001
$ws.Range("H86").Value = 1188
$ws.Range("I86").Value = 454
$ws.Range("K86").Value = 1362
$ws.Range("M86").Value = -176
$ws.Range("H89").Value = 1188
$ws.Range("I89").Value = 454
$ws.Range("K89").Value = 4086
$ws.Range("M89").Value = 1842
$ws.Range("H94").Value = 8010.2
$ws.Range("I94").Value = 1512
$ws.Range("J94").Value = 12342.333
$ws.Range("K94").Value = 4536
$ws.Range("L94").Value = 37026.999
$ws.Range("M94").Value = -3860
$ws.Range("N94").Value = -38378.999
$ws.Range("H107").Value = 1531.6923
$ws.Range("I107").Value = 1217.0555
$ws.Range("J107").Value = 2239.625
$ws.Range("K107").Value = 3651.1665
$ws.Range("L107").Value = 6718.875
$ws.Range("M107").Value = -1731.1665
$ws.Range("N107").Value = -10558.875
$ws.Range("H113").Value = 1130.6428
$ws.Range("I113").Value = 711
$ws.Range("J113").Value = 1445.375
$ws.Range("K113").Value = 2133
$ws.Range("L113").Value = 4336.125
$ws.Range("M113").Value = 37
$ws.Range("N113").Value = -8676.125
$ws.Range("H118").Value = 2724.2727
$ws.Range("I118").Value = 1994.5
$ws.Range("J118").Value = 3141.2856
$ws.Range("K118").Value = 5983.5
$ws.Range("L118").Value = 9423.856800000001
$ws.Range("M118").Value = -4740.5
$ws.Range("N118").Value = -11909.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5720.8
$ws.Range("I80").Value = 4316.0527
$ws.Range("J80").Value = 8147.1816
$ws.Range("K80").Value = 4316.0527
$ws.Range("L80").Value = 8147.1816
$ws.Range("M80").Value = -3318.0527
$ws.Range("N80").Value = -10143.1816
$ws.Range("H83").Value = 5720.8
$ws.Range("I83").Value = 4316.0527
$ws.Range("J83").Value = 8147.1816
$ws.Range("K83").Value = 21580.2635
$ws.Range("L83").Value = 40735.908
$ws.Range("M83").Value = -16588.2635
$ws.Range("N83").Value = -50719.908
$ws.Range("H102").Value = 3054.4119
$ws.Range("I102").Value = 2431.875
$ws.Range("K102").Value = 2431.875
$ws.Range("M102").Value = -809.875
$ws.Range("H128").Value = 67656
$ws.Range("J128").Value = 67656
$ws.Range("L128").Value = 67656
$ws.Range("N128").Value = -77616

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2987.5
$ws.Range("I7").Value = 2875
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 2875
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -2763
$ws.Range("N7").Value = -4224
$ws.Range("H8").Value = 100000
$ws.Range("J8").Value = 100000
$ws.Range("L8").Value = 100000
$ws.Range("N8").Value = -100280
$ws.Range("H16").Value = 3121
$ws.Range("I16").Value = 3121
$ws.Range("K16").Value = 3121
$ws.Range("M16").Value = -2951
$ws.Range("H46").Value = 5157.8887
$ws.Range("I46").Value = 1209.5
$ws.Range("K46").Value = 1209.5
$ws.Range("M46").Value = -1021.5
$ws.Range("H68").Value = 10299.8
$ws.Range("I68").Value = 7498
$ws.Range("K68").Value = 7498
$ws.Range("M68").Value = -6749
$ws.Range("H71").Value = 10299.8
$ws.Range("I71").Value = 7498
$ws.Range("K71").Value = 37490
$ws.Range("M71").Value = -33746
$ws.Range("H80").Value = 55000
$ws.Range("J80").Value = 55000
$ws.Range("L80").Value = 55000
$ws.Range("N80").Value = -57246
$ws.Range("H83").Value = 55000
$ws.Range("J83").Value = 55000
$ws.Range("L83").Value = 165000
$ws.Range("N83").Value = -176232
$ws.Range("H126").Value = 2987.5
$ws.Range("I126").Value = 2875
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 8625
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -6155
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6206.3335
$ws.Range("J62").Value = 6375.125
$ws.Range("L62").Value = 6375.125
$ws.Range("N62").Value = -7623.125
$ws.Range("H65").Value = 6206.3335
$ws.Range("J65").Value = 6375.125
$ws.Range("L65").Value = 31875.625
$ws.Range("N65").Value = -38115.625
